# Actualizacion de Utilidades y Contingencia Automatica
#
# Adds two new columns ("porcentaje_utilidades" and
# "porcentaje_contingencia") to the existing "Tabla1" table on the
# active sheet, fills in the "porcentaje_contingencia" values for the
# two data rows, and leaves the selection on D3 (matching the target
# workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook has a single table ("Tabla1") with columns:
#   id | nombre | proyecto_relacionado
$tbl = $ws.ListObjects.Item(1)

# Grow the table by two columns (id/nombre/proyecto_relacionado -> + D + E).
$tbl.ListColumns.Add() | Out-Null
$tbl.ListColumns.Add() | Out-Null

# Name the new header cells (this is what actually syncs the ListColumn
# name in this engine).
$ws.Range("D1").Value = "porcentaje_utilidades"
$ws.Range("E1").Value = "porcentaje_contingencia"

# Fill in the contingency percentage for the two existing data rows.
# (porcentaje_utilidades is left blank, matching the source data.)
$ws.Range("E2").Value = 13
$ws.Range("E3").Value = 13

# Match the final selection state recorded in the workbook.
$ws.Range("D3").Select()
